$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.959.18"
$ws.Range("E2").Value = "  -4.82%  "
$ws.Range("D3").Value = "'3.077.83"
$ws.Range("E3").Value = "  -4.76%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'538.85"
$ws.Range("E5").Value = "  -7.04%  "
$ws.Range("D6").Value = "'132.39"
$ws.Range("E6").Value = "  -12.83%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'3.077.14"
$ws.Range("E8").Value = "  -4.51%  "
$ws.Range("D9").Value = "'0.487"
$ws.Range("E9").Value = "  -5.07%  "
$ws.Range("D10").Value = "'0.154"
$ws.Range("E10").Value = "  -5.75%  "
$ws.Range("D11").Value = "'6.13"
$ws.Range("E11").Value = "  -13.21%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -6.28%  "
$ws.Range("D13").Value = "'0.0000225"
$ws.Range("E13").Value = "  -3.74%  "
$ws.Range("D14").Value = "'34.29"
$ws.Range("E14").Value = "  -10.76%  "
$ws.Range("D15").Value = "'3.537.55"
$ws.Range("E15").Value = "  -5.93%  "
$ws.Range("D16").Value = "'62.922.47"
$ws.Range("E16").Value = "  -5.04%  "
$ws.Range("D17").Value = "'0.110"
$ws.Range("E17").Value = "  -3.65%  "
$ws.Range("D18").Value = "'3.075.58"
$ws.Range("E18").Value = "  -4.68%  "
$ws.Range("D19").Value = "'6.57"
$ws.Range("E19").Value = "  -8.06%  "
$ws.Range("D20").Value = "'481.31"
$ws.Range("E20").Value = "  -11.19%  "
$ws.Range("D21").Value = "'13.23"
$ws.Range("E21").Value = "  -9.55%  "
$ws.Range("D22").Value = "'0.700"
$ws.Range("E22").Value = "  -5.99%  "
$ws.Range("D23").Value = "'7.15"
$ws.Range("E23").Value = "  -8.29%  "
$ws.Range("D24").Value = "'78.28"
$ws.Range("E24").Value = "  -3.65%  "
$ws.Range("D25").Value = "'11.99"
$ws.Range("E25").Value = "  -11.70%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D27").Value = "'2.68"
$ws.Range("E27").Value = "  -9.70%  "
$ws.Range("D28").Value = "'8.10"
$ws.Range("E28").Value = "  -13.91%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "'25.88"
$ws.Range("E30").Value = "  -6.44%  "
$ws.Range("E31").Value = "  -16.78%  "
$ws.Range("E32").Value = "  -6.98%  "
$ws.Range("D33").Value = "'58.73"
$ws.Range("E33").Value = "  +7.35%  "
$ws.Range("D34").Value = "'2.40"
$ws.Range("E34").Value = "  -13.18%  "
$ws.Range("D35").Value = "'5.93"
$ws.Range("E35").Value = "  -7.13%  "
$ws.Range("D36").Value = "'5.16"
$ws.Range("E36").Value = "  -8.60%  "
$ws.Range("D37").Value = "'462.47"
$ws.Range("E37").Value = "  -18.22%  "
$ws.Range("D38").Value = "'3.115.44"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("D39").Value = "'0.0388"
$ws.Range("E39").Value = "  -15.63%  "
$ws.Range("D40").Value = "'0.0785"
$ws.Range("E40").Value = "  -9.04%  "
$ws.Range("D41").Value = "'0.113"
$ws.Range("E41").Value = "  -12.61%  "
$ws.Range("D42").Value = "'8.01"
$ws.Range("E42").Value = "  -7.01%  "
$ws.Range("D43").Value = "'2.50"
$ws.Range("E43").Value = "  -14.11%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'0.248"
$ws.Range("E45").Value = "  -13.15%  "
$ws.Range("D46").Value = "'2.00"
$ws.Range("E46").Value = "  -14.89%  "
$ws.Range("D47").Value = "'24.18"
$ws.Range("E47").Value = "  -8.94%  "
$ws.Range("D48").Value = "'117.13"
$ws.Range("E48").Value = "  -6.25%  "
$ws.Range("E49").Value = "  -5.31%  "
$ws.Range("D50").Value = "'0.0₃0507"
$ws.Range("E50").Value = "  -8.88%  "
$ws.Range("D51").Value = "'1.98"
$ws.Range("E51").Value = "  -10.47%  "
